# Add three new "arrow label" columns (A_arrow, B_arrow, C_arrow) and three
# matching axis-title columns ("[ Ca(%) ]SiO3", "[ Mg(%) ]SiO3", "[ Fe(%) ]SiO3")
# to the "axes" worksheet, inserted right before the existing "Title" /
# "Pyroxene Classification Diagram" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("axes")

# Insert three new blank columns at D:F - this pushes the former D column
# (Title / Pyroxene Classification Diagram) to G, and extends the used
# range from A1:D2 to A1:G2.
$ws.Columns("D:F").Insert()

# Row 1 (header) - new arrow-label headers.
$ws.Range("D1").Value = "A_arrow"
$ws.Range("E1").Value = "B_arrow"
$ws.Range("F1").Value = "C_arrow"

# Row 2 (values) - new axis-title strings.
$ws.Range("D2").Value = "[ Ca(%) ]SiO3"
$ws.Range("E2").Value = "[ Mg(%) ]SiO3"
$ws.Range("F2").Value = "[ Fe(%) ]SiO3"

# Match the width of the new columns to columns A:C.
$ws.Columns("D:F").ColumnWidth = $ws.Columns("A").ColumnWidth

# Row 2 previously had a manually bumped height (15.6); restore it to the
# sheet's default auto height now that it's an ordinary text row again.
$ws.Rows(2).AutoFit()

# Leave the selection on the last newly-added cell, as in the saved file.
$ws.Range("F2").Select() | Out-Null

Write-Output "Inserted A_arrow/B_arrow/C_arrow and Ca/Mg/Fe axis-title columns on 'axes'."
